# Vendor test-data rotation: bump each sheet's D2:E6 "next batch" values,
# move the active sheet/selection, matching the target revision.

$wb = $excel.ActiveWorkbook

# --- addVendor: AT-86..90 / Auto_Vendor 86..90  ->  AT-111..115 / Auto_Vendor 111..115
$wsAdd = $wb.Worksheets.Item("addVendor")
$wsAdd.Range("D2").Value = "AT-111"
$wsAdd.Range("E2").Value = "Auto_Vendor 111"
$wsAdd.Range("D3").Value = "AT-112"
$wsAdd.Range("E3").Value = "Auto_Vendor 112"
$wsAdd.Range("D4").Value = "AT-113"
$wsAdd.Range("E4").Value = "Auto_Vendor 113"
$wsAdd.Range("D5").Value = "AT-114"
$wsAdd.Range("E5").Value = "Auto_Vendor 114"
$wsAdd.Range("D6").Value = "AT-115"
$wsAdd.Range("E6").Value = "Auto_Vendor 115"
$wsAdd.Range("D12").Select()

# --- deleteVendor: AT_DEL_116..120 / Delete_Vendor_116..120 -> AT_DEL_145..149 / Delete_Vendor_145..149
$wsDel = $wb.Worksheets.Item("deleteVendor")
$wsDel.Range("D2").Value = "AT_DEL_145"
$wsDel.Range("E2").Value = "Delete_Vendor_145"
$wsDel.Range("D3").Value = "AT_DEL_146"
$wsDel.Range("E3").Value = "Delete_Vendor_146"
$wsDel.Range("D4").Value = "AT_DEL_147"
$wsDel.Range("E4").Value = "Delete_Vendor_147"
$wsDel.Range("D5").Value = "AT_DEL_148"
$wsDel.Range("E5").Value = "Delete_Vendor_148"
$wsDel.Range("D6").Value = "AT_DEL_149"
$wsDel.Range("E6").Value = "Delete_Vendor_149"
$wsDel.Range("B28").Select()

# --- syncVendor: TE-VE-IN-80..84 / Del_Vendor 80..84 -> TE-VE-IN-105..109 / Del_Vendor 105..109
$wsSync = $wb.Worksheets.Item("syncVendor")
$wsSync.Range("D2").Value = "TE-VE-IN-105"
$wsSync.Range("E2").Value = "Del_Vendor 105"
$wsSync.Range("D3").Value = "TE-VE-IN-106"
$wsSync.Range("E3").Value = "Del_Vendor 106"
$wsSync.Range("D4").Value = "TE-VE-IN-107"
$wsSync.Range("E4").Value = "Del_Vendor 107"
$wsSync.Range("D5").Value = "TE-VE-IN-108"
$wsSync.Range("E5").Value = "Del_Vendor 108"
$wsSync.Range("D6").Value = "TE-VE-IN-109"
$wsSync.Range("E6").Value = "Del_Vendor 109"
$wsSync.Range("B26").Select()

# --- editVendor: AT_EDT-91..95 / Auto_Vendor_edit_91..95 -> AT_EDT-121..125 / Auto_Vendor_edit_121..125
# This sheet becomes the active tab, with the new selection E2:E6.
$wsEdit = $wb.Worksheets.Item("editVendor")
$wsEdit.Range("D2").Value = "AT_EDT-121"
$wsEdit.Range("E2").Value = "Auto_Vendor_edit_121"
$wsEdit.Range("D3").Value = "AT_EDT-122"
$wsEdit.Range("E3").Value = "Auto_Vendor_edit_122"
$wsEdit.Range("D4").Value = "AT_EDT-123"
$wsEdit.Range("E4").Value = "Auto_Vendor_edit_123"
$wsEdit.Range("D5").Value = "AT_EDT-124"
$wsEdit.Range("E5").Value = "Auto_Vendor_edit_124"
$wsEdit.Range("D6").Value = "AT_EDT-125"
$wsEdit.Range("E6").Value = "Auto_Vendor_edit_125"
$wsEdit.Activate()
$wsEdit.Range("E2:E6").Select()
